$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Limiti": update the monthly limit (B) and remaining (C) columns.
# ---------------------------------------------------------------------------
$limiti = $wb.Worksheets.Item("Limiti")

$limitiValues = @(23, 30, 100, 23, 5, 40, 2, 5, 2)
for ($i = 0; $i -lt $limitiValues.Length; $i++) {
    $row = $i + 2
    $limiti.Cells.Item($row, 2).Value = $limitiValues[$i]
    $limiti.Cells.Item($row, 3).Value = $limitiValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Izdevumi": refresh existing rows 2-8 and append new rows 9-25.
# ---------------------------------------------------------------------------
$izdevumi = $wb.Worksheets.Item("Izdevumi")

$rows = @(
    @{ A = "2025-05-09 18:55:00"; B = "pārtika"; C = 17 },
    @{ A = "2025-05-09 18:55:11"; B = "ēšana ārpus mājas (restorāni/fast food/kafejnīcas)"; C = 23.32 },
    @{ A = "2025-05-09 18:55:39"; B = "mājas izdevumi (komunālie + īre / nekustamā īpašuma nodoklis)"; C = 111.5 },
    @{ A = "2025-05-09 18:55:45"; B = "hobiji"; C = 12.3 },
    @{ A = "2025-05-09 18:55:48"; B = "mājdzīvnieki"; C = 0 },
    @{ A = "2025-05-09 18:55:56"; B = "apģērbs"; C = 13.99 },
    @{ A = "2025-05-09 18:55:58"; B = "higēnas preces"; C = 0 },
    @{ A = "2025-05-09 18:56:01"; B = "medicīniskie izdevumi"; C = 0 },
    @{ A = "2025-05-09 18:56:05"; B = "transports"; C = 13 },
    @{ A = "2025-05-09 18:56:09"; B = "izklaide (kino, teātris, klubs)"; C = 17.5 },
    @{ A = "2025-05-09 18:56:12"; B = "abonementi"; C = 19.99 },
    @{ A = "2025-05-09 18:56:14"; B = "dāvanas"; C = 0 },
    @{ A = "2025-05-09 18:56:28"; B = "pārtika"; C = 2 },
    @{ A = "2025-05-09 18:56:30"; B = "ēšana ārpus mājas (restorāni/fast food/kafejnīcas)"; C = 0 },
    @{ A = "2025-05-09 18:56:31"; B = "mājas izdevumi (komunālie + īre / nekustamā īpašuma nodoklis)"; C = 0 },
    @{ A = "2025-05-09 18:56:31"; B = "hobiji"; C = 0 },
    @{ A = "2025-05-09 18:56:32"; B = "mājdzīvnieki"; C = 0 },
    @{ A = "2025-05-09 18:56:32"; B = "apģērbs"; C = 0 },
    @{ A = "2025-05-09 18:56:33"; B = "higēnas preces"; C = 0 },
    @{ A = "2025-05-09 18:56:34"; B = "medicīniskie izdevumi"; C = 0 },
    @{ A = "2025-05-09 18:56:35"; B = "transports"; C = 0 },
    @{ A = "2025-05-09 18:56:36"; B = "izklaide (kino, teātris, klubs)"; C = 0 },
    @{ A = "2025-05-09 18:56:36"; B = "abonementi"; C = 0 },
    @{ A = "2025-05-09 18:56:37"; B = "dāvanas"; C = 0 }
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $i + 2
    $izdevumi.Cells.Item($row, 1).Value = $rows[$i].A
    $izdevumi.Cells.Item($row, 2).Value = $rows[$i].B
    $izdevumi.Cells.Item($row, 3).Value = $rows[$i].C
}
